# Update the valid registration scenario data on the Register_valid sheet.
$wb = $excel.ActiveWorkbook

$wsInvalid = $wb.Worksheets.Item("Register_invalid")
$wsValid = $wb.Worksheets.Item("Register_valid")

# Update shared values used by the valid-entry scenario row.
$wsValid.Range("B2").Value = "mita@671"
$wsValid.Range("C2").Value = "qwerew123"
$wsValid.Range("D2").Value = "qwerew123"

# Update the active selections to match the authored state.
$wsInvalid.Range("F12").Select() | Out-Null
$wsValid.Range("D2").Select() | Out-Null
